# Fixed camera bug. Free camera now seems to do what I want it to!
# Delete the entire row 3 (which contained the duplicate/incorrect
# "Investigate camera issues" task), shifting all rows below up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet carries 3 legacy cell-notes (B5, B13, B19). They are anchored
# to absolute cells and won't automatically follow the upward shift caused
# by deleting row 3, so capture their text first and re-home them after
# the row is gone.
$note1 = $ws.Range("B5").Comment.Text()
$note2 = $ws.Range("B13").Comment.Text()
$note3 = $ws.Range("B19").Comment.Text()

$ws.Rows.Item(3).Delete()

# Restore the selection to mirror what Excel leaves selected after an
# entire-row delete: the whole of the (new) row 3.
$ws.Range("A3:XFD3").Select()

# Move the notes up by one row (B5->B4, B13->B12, B19->B18) to follow
# the cells they originally annotated.
$ws.Range("B5").Comment.Delete()
$ws.Range("B13").Comment.Delete()
$ws.Range("B19").Comment.Delete()

$ws.Range("B4").AddComment($note1)
$ws.Range("B12").AddComment($note2)
$ws.Range("B18").AddComment($note3)
